$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.245.08'
$ws.Range("E2").Value = '  +3.22%  '
$ws.Range("D3").Value = '3.072.79'
$ws.Range("E3").Value = '  +5.69%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '514.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.07'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.42%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  +3.85%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.22'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.91%  '
$ws.Range("E10").Value = '  +4.59%  '
$ws.Range("E11").Value = '  +6.78%  '
$ws.Range("D12").Value = '3.597.57'
$ws.Range("E12").Value = '  +6.00%  '
$ws.Range("E13").Value = '  +2.84%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.48'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000165'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.13%  '
$ws.Range("D16").Value = '57.308.01'
$ws.Range("E16").Value = '  +3.39%  '
$ws.Range("D17").Value = '3.078.29'
$ws.Range("E17").Value = '  +6.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.93'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.03'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.45%  '
$ws.Range("E20").Value = '  +7.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '336.52'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.33%  '
$ws.Range("E22").Value = '  +0.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.501'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.30'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.21%  '
$ws.Range("E25").Value = '  +6.77%  '
$ws.Range("B26").Value = 'Binance-PegBSC-USD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.20%  '
$ws.Range("B27").Value = 'PEPE'
$ws.Range("C27").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D27").Value = '0.0₃0947'
$ws.Range("E27").Value = '  +12.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.47'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.07'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.70%  '
$ws.Range("E30").Value = '  +3.60%  '
$ws.Range("E31").Value = '  +5.70%  '
$ws.Range("E32").Value = '  +6.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '154.16'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.56%  '
$ws.Range("E34").Value = '  +4.35%  '
$ws.Range("E35").Value = '  +5.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '26.23'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.37%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.24'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0672'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.89%  '
$ws.Range("D39").Value = '3.114.20'
$ws.Range("E39").Value = '  +6.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.27%  '
$ws.Range("E41").Value = '  +6.01%  '
$ws.Range("E42").Value = '  +4.82%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").Value = '2.242.16'
$ws.Range("E44").Value = '  +7.50%  '
$ws.Range("E45").Value = '  +10.03%  '
$ws.Range("E46").Value = '  +5.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.949'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.99%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.01'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.09%  '
$ws.Range("E49").Value = '  -0.23%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0867'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.85%  '
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.74'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.46%  '
